$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Turn off iterative calculation (workbook-level setting)
$excel.Iterative = $false

# Add new row 51 data
$ws.Range("A51").Value = "XXX56e1b-bc07-41cd-bad4-a5b51b6287da"
$ws.Range("B51").Value = "NEW DUMMY"
$ws.Range("C51").Value = "kg/m3"
$ws.Range("D51").Value = "NEW DUMMY COMMENT"
$ws.Range("E51").Value = 666
$ws.Range("F51").Value = 666
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 1
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = 1
$ws.Range("K51").Value = 1
$ws.Range("L51").Value = "x"
$ws.Range("M51").Value = "x"
$ws.Range("N51").Value = "x"
$ws.Range("O51").Value = "x"

# Style the new row to mirror existing ones
$ws.Range("L51").Style = $ws.Range("L50").Style
$ws.Range("O51").Style = $ws.Range("O49").Style

# Selection / view state changes
$ws.Range("A51:XFD51").Select()
$ws.Application.ActiveWindow.ScrollRow = 36
